$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
try {
  $shape.TextEffect.PresetTextEffect = 1
  Write-Output "set ok"
} catch {
  Write-Output "ERR: $_"
}
